$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 283.875
$ws.Range("I33").Value = 274.68
$ws.Range("J33").Value = 316.7143
$ws.Range("K33").Value = 274.68
$ws.Range("L33").Value = 316.7143
$ws.Range("M33").Value = -45.68000000000001
$ws.Range("N33").Value = -774.7143
$ws.Range("H132").Value = 593096.0600000001
$ws.Range("I132").Value = 2534.9155
$ws.Range("K132").Value = 7604.7465
$ws.Range("M132").Value = -5074.7465
$ws.Range("H137").Value = 2441855.2
$ws.Range("I137").Value = 4350152
$ws.Range("K137").Value = 13050456
$ws.Range("M137").Value = -13047906
$ws.Range("H138").Value = 2820390
$ws.Range("I138").Value = 3102.9285
$ws.Range("J138").Value = 3512355.2
$ws.Range("K138").Value = 9308.7855
$ws.Range("L138").Value = 10537065.6
$ws.Range("M138").Value = -4168.7855
$ws.Range("N138").Value = -10547345.6

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22080.42
$ws.Range("I32").Value = 12963.058
$ws.Range("J32").Value = 42373.902
$ws.Range("K32").Value = 12963.058
$ws.Range("L32").Value = 42373.902
$ws.Range("M32").Value = -12676.058
$ws.Range("N32").Value = -42947.902
$ws.Range("H61").Value = 33401782
$ws.Range("I61").Value = 40041772
$ws.Range("K61").Value = 40041772
$ws.Range("M61").Value = -40041560
$ws.Range("H132").Value = 15692801
$ws.Range("I132").Value = 22775150
$ws.Range("J132").Value = 111634
$ws.Range("K132").Value = 68325450
$ws.Range("L132").Value = 334902
$ws.Range("M132").Value = -68322920
$ws.Range("N132").Value = -339962
$ws.Range("H136").Value = 33401782
$ws.Range("I136").Value = 40041772
$ws.Range("K136").Value = 120125316
$ws.Range("M136").Value = -120122766

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 445.05884
$ws.Range("I94").Value = 379.125
$ws.Range("K94").Value = 379.125
$ws.Range("M94").Value = 71.875

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 227.6
$ws.Range("I7").Value = 119
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 119
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -6
$ws.Range("N7").Value = -526
$ws.Range("H10").Value = 327.6
$ws.Range("I10").Value = 327.6
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 327.6
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = ""
$ws.Range("N10").Value = -188.6
$ws.Range("H31").Value = 41395.09
$ws.Range("I31").Value = 35480.5
$ws.Range("J31").Value = 48219.617
$ws.Range("K31").Value = 35480.5
$ws.Range("L31").Value = 48219.617
$ws.Range("M31").Value = -35185.5
$ws.Range("N31").Value = -48809.617
$ws.Range("H34").Value = 41395.09
$ws.Range("I34").Value = 35480.5
$ws.Range("J34").Value = 48219.617
$ws.Range("K34").Value = 35480.5
$ws.Range("L34").Value = 48219.617
$ws.Range("M34").Value = -35278.5
$ws.Range("N34").Value = -48623.617
$ws.Range("H132").Value = 33524.562
$ws.Range("I132").Value = 2013
$ws.Range("J132").Value = 86043.836
$ws.Range("K132").Value = 6039
$ws.Range("L132").Value = 258131.508
$ws.Range("M132").Value = -3509
$ws.Range("N132").Value = -263191.508

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1340.5
$ws.Range("I69").Value = 674.8570999999999
$ws.Range("J69").Value = 6000
$ws.Range("K69").Value = 2024.5713
$ws.Range("L69").Value = 18000
$ws.Range("M69").Value = -1213.5713
$ws.Range("N69").Value = -19622
$ws.Range("H72").Value = 1340.5
$ws.Range("I72").Value = 674.8570999999999
$ws.Range("J72").Value = 6000
$ws.Range("K72").Value = 6073.7139
$ws.Range("L72").Value = 54000
$ws.Range("M72").Value = -2017.7139
$ws.Range("N72").Value = -62112
$ws.Range("H131").Value = 910.54
$ws.Range("I131").Value = 393.33334
$ws.Range("J131").Value = 943.55316
$ws.Range("K131").Value = 1180.00002
$ws.Range("L131").Value = 2830.65948
$ws.Range("M131").Value = 3859.99998
$ws.Range("N131").Value = -12910.65948

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 173.6
$ws.Range("I2").Value = 122
$ws.Range("J2").Value = 251
$ws.Range("K2").Value = 122
$ws.Range("L2").Value = 251
$ws.Range("M2").Value = -9
$ws.Range("N2").Value = -477
$ws.Range("H3").Value = 1225.8
$ws.Range("I3").Value = 1650
$ws.Range("J3").Value = 943
$ws.Range("K3").Value = 1650
$ws.Range("L3").Value = 943
$ws.Range("M3").Value = -1534
$ws.Range("N3").Value = -1175
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("H18").Value = 6000
$ws.Range("J18").Value = 6000
$ws.Range("L18").Value = 6000
$ws.Range("N18").Value = -6586
$ws.Range("H43").Value = 15000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = ""
$ws.Range("M43").Value = 15000
$ws.Range("N43").Value = -15302
$ws.Range("H46").Value = 25554.1
$ws.Range("J46").Value = 36428.57
$ws.Range("L46").Value = 36428.57
$ws.Range("N46").Value = -36740.57
$ws.Range("H80").Value = 4018.468
$ws.Range("I80").Value = 3983.3333
$ws.Range("J80").Value = 4023.6099
$ws.Range("K80").Value = 3983.3333
$ws.Range("L80").Value = 4023.6099
$ws.Range("M80").Value = -2985.3333
$ws.Range("N80").Value = -6019.609899999999
$ws.Range("H83").Value = 4018.468
$ws.Range("I83").Value = 3983.3333
$ws.Range("J83").Value = 4023.6099
$ws.Range("K83").Value = 19916.6665
$ws.Range("L83").Value = 20118.0495
$ws.Range("M83").Value = -14924.6665
$ws.Range("N83").Value = -30102.0495
$ws.Range("H93").Value = 18531.375
$ws.Range("I93").Value = 18000
$ws.Range("J93").Value = 20125.5
$ws.Range("K93").Value = 18000
$ws.Range("L93").Value = 20125.5
$ws.Range("M93").Value = -16128
$ws.Range("N93").Value = -23869.5
$ws.Range("H113").Value = 2269.8333
$ws.Range("I113").Value = 1554
$ws.Range("J113").Value = 2875.5386
$ws.Range("K113").Value = 1554
$ws.Range("L113").Value = 2875.5386
$ws.Range("M113").Value = 616
$ws.Range("N113").Value = -7215.5386

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2600
$ws.Range("I2").Value = 1666.6666
$ws.Range("K2").Value = 1666.6666
$ws.Range("M2").Value = -1554.6666
$ws.Range("H16").Value = 2675.5
$ws.Range("I16").Value = 2766.6667
$ws.Range("J16").Value = 2402
$ws.Range("K16").Value = 2766.6667
$ws.Range("L16").Value = 2402
$ws.Range("M16").Value = -2596.6667
$ws.Range("N16").Value = -2742
$ws.Range("H22").Value = 909.6
$ws.Range("I22").Value = 815.2
$ws.Range("K22").Value = 815.2
$ws.Range("M22").Value = -520.2
$ws.Range("H27").Value = 909.6
$ws.Range("I27").Value = 815.2
$ws.Range("K27").Value = 815.2
$ws.Range("M27").Value = -708.2
$ws.Range("H93").Value = 1296.6666
$ws.Range("I93").Value = 1296.6666
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1296.6666
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = -48.66660000000002

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 554.88
$ws.Range("I107").Value = 427.7619
$ws.Range("K107").Value = 1283.2857
$ws.Range("M107").Value = 636.7142999999999
$ws.Range("H113").Value = 800.3929000000001
$ws.Range("I113").Value = 990.3570999999999
$ws.Range("J113").Value = 610.4286
$ws.Range("K113").Value = 2971.0713
$ws.Range("L113").Value = 1831.2858
$ws.Range("M113").Value = -801.0712999999996
$ws.Range("N113").Value = -6171.2858
$ws.Range("H132").Value = 36933.34
$ws.Range("I132").Value = 25839.45
$ws.Range("J132").Value = 64668.062
$ws.Range("K132").Value = 77518.35000000001
$ws.Range("L132").Value = 194004.186
$ws.Range("M132").Value = -74988.35000000001
$ws.Range("N132").Value = -199064.186
$ws.Range("H136").Value = 38806.836
$ws.Range("I136").Value = 26286.65
$ws.Range("J136").Value = 72194
$ws.Range("K136").Value = 78859.95000000001
$ws.Range("L136").Value = 216582
$ws.Range("M136").Value = -76309.95000000001
$ws.Range("N136").Value = -221682
